$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("ALC")
$ws2 = $wb.Worksheets.Item("ARM")
$ws3 = $wb.Worksheets.Item("BSM")
$ws4 = $wb.Worksheets.Item("CRP")
$ws5 = $wb.Worksheets.Item("CUL")
$ws6 = $wb.Worksheets.Item("GSM")
$ws7 = $wb.Worksheets.Item("LTW")
$ws8 = $wb.Worksheets.Item("WVR")

# --- Sheet ALC ---
# row 49 (hunk @@ -3069,25 +3069,25 @@)
$ws1.Range("H49").Value = 352.25
$ws1.Range("I49").Value = 290
$ws1.Range("J49").Value = 373
$ws1.Range("K49").Value = 870
$ws1.Range("L49").Value = 1119
$ws1.Range("M49").Value = -734
$ws1.Range("N49").Value = -1391

# row 92 (hunk @@ -5236,25 +5236,25 @@)
$ws1.Range("H92").Value = 8548291
$ws1.Range("I92").Value = 13334199
$ws1.Range("J92").Value = 2027.4286
$ws1.Range("K92").Value = 13334199
$ws1.Range("L92").Value = 2027.4286
$ws1.Range("M92").Value = -13332951
$ws1.Range("N92").Value = -4523.4286

# row 135 (hunk @@ -7394,25 +7394,25 @@)
$ws1.Range("H135").Value = 946.57574
$ws1.Range("I135").Value = 524.8214
$ws1.Range("J135").Value = 3308.4
$ws1.Range("K135").Value = 4723.3926
$ws1.Range("L135").Value = 29775.6
$ws1.Range("M135").Value = -2188.3926
$ws1.Range("N135").Value = -34845.60000000001

# --- Sheet ARM ---
# row 2 (hunk @@ -7843,22 +7843,22 @@)
$ws2.Range("H2").Value = 653
$ws2.Range("I2").Value = 552.25
$ws2.Range("K2").Value = 552.25
$ws2.Range("M2").Value = -439.25

# row 45 (hunk @@ -9938,25 +9938,25 @@)
$ws2.Range("H45").Value = 1884.1904
$ws2.Range("I45").Value = 1446.6666
$ws2.Range("J45").Value = 2978
$ws2.Range("K45").Value = 1446.6666
$ws2.Range("L45").Value = 2978
$ws2.Range("M45").Value = -1069.6666
$ws2.Range("N45").Value = -3732

# row 88 (hunk @@ -12030,25 +12030,22 @@)
$ws2.Range("H88").Value = 3000
$ws2.Range("I88").Value = 0
$ws2.Range("J88").Value = 3000
$ws2.Range("K88").Value = 0
$ws2.Range("L88").Value = 3000
$ws2.Range("M88").ClearContents()
$ws2.Range("N88").Value = -3812

# row 91 (hunk @@ -12180,25 +12177,22 @@)
$ws2.Range("H91").Value = 3000
$ws2.Range("I91").Value = 0
$ws2.Range("J91").Value = 3000
$ws2.Range("K91").Value = 0
$ws2.Range("L91").Value = 3000
$ws2.Range("M91").ClearContents()
$ws2.Range("N91").Value = -5808

# row 97 (hunk @@ -12477,25 +12471,25 @@)
$ws2.Range("H97").Value = 1172.5
$ws2.Range("I97").Value = 1006.36365
$ws2.Range("J97").Value = 3000
$ws2.Range("K97").Value = 1006.36365
$ws2.Range("L97").Value = 3000
$ws2.Range("M97").Value = -510.36365
$ws2.Range("N97").Value = -3992

# row 116 (hunk @@ -13408,22 +13402,22 @@)
$ws2.Range("H116").Value = 653
$ws2.Range("I116").Value = 552.25
$ws2.Range("K116").Value = 552.25
$ws2.Range("M116").Value = 1741.75

# row 132 (hunk @@ -14189,25 +14183,25 @@)
$ws2.Range("H132").Value = 3822.5938
$ws2.Range("I132").Value = 3916.25
$ws2.Range("J132").Value = 3666.5
$ws2.Range("K132").Value = 11748.75
$ws2.Range("L132").Value = 10999.5
$ws2.Range("M132").Value = -9218.75
$ws2.Range("N132").Value = -16059.5

# --- Sheet BSM ---
# row 3 (hunk @@ -14831,22 +14825,22 @@)
$ws3.Range("H3").Value = 653
$ws3.Range("I3").Value = 552.25
$ws3.Range("K3").Value = 552.25
$ws3.Range("M3").Value = -438.25

# row 86 (hunk @@ -18910,25 +18904,25 @@)
$ws3.Range("H86").Value = 2768.6428
$ws3.Range("I86").Value = 2384.25
$ws3.Range("J86").Value = 3281.1667
$ws3.Range("K86").Value = 2384.25
$ws3.Range("L86").Value = 3281.1667
$ws3.Range("M86").Value = -1261.25
$ws3.Range("N86").Value = -5527.1667

# row 89 (hunk @@ -19060,25 +19054,25 @@)
$ws3.Range("H89").Value = 2768.6428
$ws3.Range("I89").Value = 2384.25
$ws3.Range("J89").Value = 3281.1667
$ws3.Range("K89").Value = 11921.25
$ws3.Range("L89").Value = 16405.8335
$ws3.Range("M89").Value = -6305.25
$ws3.Range("N89").Value = -27637.8335

# row 94 (hunk @@ -19308,25 +19302,25 @@)
$ws3.Range("H94").Value = 310
$ws3.Range("I94").Value = 332
$ws3.Range("J94").Value = 255
$ws3.Range("K94").Value = 332
$ws3.Range("L94").Value = 255
$ws3.Range("M94").Value = 119
$ws3.Range("N94").Value = -1157

# row 105 (hunk @@ -19856,25 +19850,22 @@)
$ws3.Range("H105").Value = 41668830
$ws3.Range("I105").Value = 41668830
$ws3.Range("J105").Value = 0
$ws3.Range("K105").Value = 41668830
$ws3.Range("L105").Value = 0
$ws3.Range("M105").Value = -41667083
$ws3.Range("N105").ClearContents()

# --- Sheet CRP ---
# row 105 (hunk @@ -26828,25 +26819,25 @@)
$ws4.Range("H105").Value = 790.7143
$ws4.Range("J105").Value = 712.2222
$ws4.Range("L105").Value = 712.2222
$ws4.Range("N105").Value = -4206.2222

# row 132 (hunk @@ -28154,25 +28145,25 @@)
$ws4.Range("H132").Value = 5052359.5
$ws4.Range("I132").Value = 1675.36
$ws4.Range("J132").Value = 20835748
$ws4.Range("K132").Value = 5026.08
$ws4.Range("L132").Value = 62507244
$ws4.Range("M132").Value = -2496.08
$ws4.Range("N132").Value = -62512304

# --- Sheet CUL ---
# row 88 (hunk @@ -33150,22 +33141,22 @@)
$ws5.Range("H88").Value = 6394.45
$ws5.Range("J88").Value = 6394.45
$ws5.Range("L88").Value = 19183.35
$ws5.Range("N88").Value = -20039.35

# row 91 (hunk @@ -33303,22 +33294,22 @@)
$ws5.Range("H91").Value = 6394.45
$ws5.Range("J91").Value = 6394.45
$ws5.Range("L91").Value = 19183.35
$ws5.Range("N91").Value = -22147.35

# row 131 (hunk @@ -35341,25 +35332,25 @@)
$ws5.Range("H131").Value = 3697.9756
$ws5.Range("I131").Value = 442.85715
$ws5.Range("J131").Value = 5385.815
$ws5.Range("K131").Value = 1328.57145
$ws5.Range("L131").Value = 16157.445
$ws5.Range("M131").Value = 3711.42855
$ws5.Range("N131").Value = -26237.445

# row 137 (hunk @@ -35653,25 +35644,25 @@)
$ws5.Range("H137").Value = 45185.42
$ws5.Range("I137").Value = 9052.532999999999
$ws5.Range("J137").Value = 94457.55
$ws5.Range("K137").Value = 27157.599
$ws5.Range("L137").Value = 283372.65
$ws5.Range("M137").Value = -22057.599
$ws5.Range("N137").Value = -293572.65

# row 140 (hunk @@ -35809,25 +35800,25 @@)
$ws5.Range("H140").Value = 1365.4584
$ws5.Range("I140").Value = 1173.3334
$ws5.Range("J140").Value = 2710.3333
$ws5.Range("K140").Value = 3520.0002
$ws5.Range("L140").Value = 8130.999899999999
$ws5.Range("M140").Value = 1659.9998
$ws5.Range("N140").Value = -18490.9999

# row 141 (hunk @@ -35861,25 +35852,25 @@)
$ws5.Range("H141").Value = 10874.737
$ws5.Range("I141").Value = 11402.857
$ws5.Range("J141").Value = 10566.667
$ws5.Range("K141").Value = 34208.571
$ws5.Range("L141").Value = 31700.001
$ws5.Range("M141").Value = -29028.571
$ws5.Range("N141").Value = -42060.001

# --- Sheet GSM ---
# row 80 (hunk @@ -39787,25 +39778,25 @@)
$ws6.Range("H80").Value = 17608768
$ws6.Range("I80").Value = 21960584
$ws6.Range("J80").Value = 201500
$ws6.Range("K80").Value = 21960584
$ws6.Range("L80").Value = 201500
$ws6.Range("M80").Value = -21959586
$ws6.Range("N80").Value = -203496

# row 83 (hunk @@ -39937,25 +39928,25 @@)
$ws6.Range("H83").Value = 17608768
$ws6.Range("I83").Value = 21960584
$ws6.Range("J83").Value = 201500
$ws6.Range("K83").Value = 109802920
$ws6.Range("L83").Value = 1007500
$ws6.Range("M83").Value = -109797928
$ws6.Range("N83").Value = -1017484

# row 94 (hunk @@ -40479,22 +40470,22 @@)
$ws6.Range("H94").Value = 95448
$ws6.Range("J94").Value = 95448
$ws6.Range("L94").Value = 95448
$ws6.Range("N94").Value = -96800

# row 102 (hunk @@ -40874,22 +40865,22 @@)
$ws6.Range("H102").Value = 1728
$ws6.Range("I102").Value = 1564.8
$ws6.Range("K102").Value = 1564.8
$ws6.Range("M102").Value = 57.20000000000005

# row 132 (hunk @@ -42332,25 +42323,25 @@)
$ws6.Range("H132").Value = 2394.319
$ws6.Range("I132").Value = 1773.7273
$ws6.Range("J132").Value = 3857.1428
$ws6.Range("K132").Value = 5321.1819
$ws6.Range("L132").Value = 11571.4284
$ws6.Range("M132").Value = -2791.1819
$ws6.Range("N132").Value = -16631.4284

# --- Sheet LTW ---
# row 132 (hunk @@ -49274,25 +49265,25 @@)
$ws7.Range("H132").Value = 3510.16
$ws7.Range("I132").Value = 2904
$ws7.Range("J132").Value = 4798.25
$ws7.Range("K132").Value = 8712
$ws7.Range("L132").Value = 14394.75
$ws7.Range("M132").Value = -6182
$ws7.Range("N132").Value = -19454.75

# row 136 (hunk @@ -49470,22 +49461,22 @@)
$ws7.Range("H136").Value = 3334800.2
$ws7.Range("I136").Value = 1251.75
$ws7.Range("K136").Value = 3755.25
$ws7.Range("M136").Value = -1205.25

# --- Sheet WVR ---
# row 81 (hunk @@ -53717,25 +53708,25 @@)
$ws8.Range("H81").Value = 3559.3333
$ws8.Range("I81").Value = 3328.4614
$ws8.Range("J81").Value = 4159.6
$ws8.Range("K81").Value = 6656.9228
$ws8.Range("L81").Value = 8319.200000000001
$ws8.Range("M81").Value = -5595.9228
$ws8.Range("N81").Value = -10441.2

# row 84 (hunk @@ -53867,25 +53858,25 @@)
$ws8.Range("H84").Value = 3559.3333
$ws8.Range("I84").Value = 3328.4614
$ws8.Range("J84").Value = 4159.6
$ws8.Range("K84").Value = 33284.614
$ws8.Range("L84").Value = 41596
$ws8.Range("M84").Value = -27980.614
$ws8.Range("N84").Value = -52204

# row 100 (hunk @@ -54642,25 +54633,25 @@)
$ws8.Range("H100").Value = 940.5
$ws8.Range("I100").Value = 880
$ws8.Range("J100").Value = 1001
$ws8.Range("K100").Value = 1760
$ws8.Range("L100").Value = 2002
$ws8.Range("M100").Value = -1219
$ws8.Range("N100").Value = -3084

# row 107 (hunk @@ -54985,25 +54976,25 @@)
$ws8.Range("H107").Value = 584.4
$ws8.Range("I107").Value = 507.27274
$ws8.Range("J107").Value = 678.6667
$ws8.Range("K107").Value = 1521.81822
$ws8.Range("L107").Value = 2036.0001
$ws8.Range("M107").Value = 398.1817799999999
$ws8.Range("N107").Value = -5876.0001
